$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (shifts existing D..K data to E..L)
$ws.Columns("D:D").Insert()

# Copy formats from column E (the old D, now shifted) into the new column D
# so the new column visually matches its neighbours (date style / number style).
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the newly-inserted column D with the latest reporting-period figures.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 480200
$ws.Range("D9").Value = 322100
$ws.Range("D10").Value = 158100
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 410600
$ws.Range("D18").Value = 69500
$ws.Range("D20").Value = -10500
$ws.Range("D21").Value = 86100
$ws.Range("D22").Value = 20000
$ws.Range("D23").Value = 39000
$ws.Range("D24").Value = 11800
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 27200
$ws.Range("D27").Value = 27200
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 10500
$ws.Range("D33").Value = 27200
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 27200
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 250100
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 57000
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 561700
$ws.Range("D48").Value = 18000
$ws.Range("D49").Value = 229300
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 12300
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 1768700
$ws.Range("D57").Value = 0
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 1080900
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 148800
$ws.Range("D62").Value = 7700
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1343400
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 195800
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 425300
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 27200
$ws.Range("D83").Value = 27100
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 95400
$ws.Range("D91").Value = -2300
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 24400
$ws.Range("D96").Value = -6400
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -32000
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 87800

# Two line items also had their prior-year figures restated alongside the new column.
$ws.Range("E89").Value = 7500
$ws.Range("F89").Value = 83000
$ws.Range("E102").Value = 47800
$ws.Range("F102").Value = -122600
